$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1490.7
$ws.Range("I2").Value = 1485.1428
$ws.Range("K2").Value = 1485.1428
$ws.Range("M2").Value = -1372.1428
$ws.Range("H19").Value = 3975.1428
$ws.Range("I19").Value = 3232.6667
$ws.Range("K19").Value = 3232.6667
$ws.Range("M19").Value = -3057.6667
$ws.Range("H32").Value = 8801.058999999999
$ws.Range("J32").Value = 5593.9165
$ws.Range("L32").Value = 5593.9165
$ws.Range("N32").Value = -6245.9165
$ws.Range("H40").Value = 3824.875
$ws.Range("J40").Value = 5966.3335
$ws.Range("L40").Value = 5966.3335
$ws.Range("N40").Value = -6316.3335
$ws.Range("H62").Value = 5634.357
$ws.Range("I62").Value = 5681.25
$ws.Range("K62").Value = 5681.25
$ws.Range("M62").Value = -5057.25
$ws.Range("H65").Value = 5634.357
$ws.Range("I65").Value = 5681.25
$ws.Range("K65").Value = 28406.25
$ws.Range("M65").Value = -25286.25
$ws.Range("H113").Value = 5598.6
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 5598.6
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -12106.6
$ws.Range("H132").Value = 17517.592
$ws.Range("I132").Value = 19421.281
$ws.Range("K132").Value = 58263.84299999999
$ws.Range("M132").Value = -55733.84299999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 734.1429000000001
$ws.Range("I2").Value = 734.1429000000001
$ws.Range("K2").Value = 734.1429000000001
$ws.Range("M2").Value = -621.1429000000001
$ws.Range("H32").Value = 23880.09
$ws.Range("I32").Value = 23880.09
$ws.Range("K32").Value = 23880.09
$ws.Range("M32").Value = -23593.09
$ws.Range("H45").Value = 3489
$ws.Range("I45").Value = 2217.375
$ws.Range("K45").Value = 2217.375
$ws.Range("M45").Value = -1840.375
$ws.Range("H61").Value = 4158.0625
$ws.Range("I61").Value = 787.8214
$ws.Range("K61").Value = 787.8214
$ws.Range("M61").Value = -575.8214
$ws.Range("H74").Value = 226642.3
$ws.Range("I74").Value = 261459.17
$ws.Range("K74").Value = 261459.17
$ws.Range("M74").Value = -260585.17
$ws.Range("H77").Value = 226642.3
$ws.Range("I77").Value = 261459.17
$ws.Range("K77").Value = 1307295.85
$ws.Range("M77").Value = -1302927.85
$ws.Range("H116").Value = 734.1429000000001
$ws.Range("I116").Value = 734.1429000000001
$ws.Range("K116").Value = 734.1429000000001
$ws.Range("M116").Value = 1559.8571
$ws.Range("H136").Value = 4158.0625
$ws.Range("I136").Value = 787.8214
$ws.Range("K136").Value = 2363.4642
$ws.Range("M136").Value = 186.5357999999997

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 734.1429000000001
$ws.Range("I3").Value = 734.1429000000001
$ws.Range("K3").Value = 734.1429000000001
$ws.Range("M3").Value = -620.1429000000001
$ws.Range("H86").Value = 6012.25
$ws.Range("I86").Value = 4979.8
$ws.Range("J86").Value = 7733
$ws.Range("K86").Value = 4979.8
$ws.Range("L86").Value = 7733
$ws.Range("M86").Value = -3856.8
$ws.Range("N86").Value = -9979
$ws.Range("H89").Value = 6012.25
$ws.Range("I89").Value = 4979.8
$ws.Range("J89").Value = 7733
$ws.Range("K89").Value = 24899
$ws.Range("L89").Value = 38665
$ws.Range("M89").Value = -19283
$ws.Range("N89").Value = -49897

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1274.4524
$ws.Range("I58").Value = 990.65717
$ws.Range("K58").Value = 990.65717
$ws.Range("M58").Value = -787.65717
$ws.Range("H134").Value = 1440.1111
$ws.Range("I134").Value = 1192.0571
$ws.Range("J134").Value = 2308.3
$ws.Range("K134").Value = 3576.1713
$ws.Range("L134").Value = 6924.900000000001
$ws.Range("M134").Value = -1041.1713
$ws.Range("N134").Value = -11994.9
$ws.Range("H136").Value = 1274.4524
$ws.Range("I136").Value = 990.65717
$ws.Range("K136").Value = 2971.97151
$ws.Range("M136").Value = -421.9715099999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 4854.95
$ws.Range("I62").Value = 2100
$ws.Range("J62").Value = 4999.9473
$ws.Range("K62").Value = 6300
$ws.Range("L62").Value = 14999.8419
$ws.Range("M62").Value = -5614
$ws.Range("N62").Value = -16371.8419
$ws.Range("H65").Value = 4854.95
$ws.Range("I65").Value = 2100
$ws.Range("J65").Value = 4999.9473
$ws.Range("K65").Value = 18900
$ws.Range("L65").Value = 44999.5257
$ws.Range("M65").Value = -15468
$ws.Range("N65").Value = -51863.5257
$ws.Range("H107").Value = 3169.5386
$ws.Range("I107").Value = 7801.5
$ws.Range("J107").Value = 1110.8889
$ws.Range("K107").Value = 23404.5
$ws.Range("L107").Value = 3332.6667
$ws.Range("M107").Value = -21484.5
$ws.Range("N107").Value = -7172.6667
$ws.Range("H131").Value = 4098.0835
$ws.Range("I131").Value = 5682.857
$ws.Range("J131").Value = 1879.4
$ws.Range("K131").Value = 17048.571
$ws.Range("L131").Value = 5638.200000000001
$ws.Range("M131").Value = -12008.571
$ws.Range("N131").Value = -15718.2
$ws.Range("H140").Value = 3085.1667
$ws.Range("I140").Value = 3085.1667
$ws.Range("K140").Value = 9255.500100000001
$ws.Range("M140").Value = -4075.500100000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H97").Value = 676.1667
$ws.Range("I97").Value = 762.38464
$ws.Range("K97").Value = 762.38464
$ws.Range("M97").Value = -266.38464
$ws.Range("H102").Value = 15144.054
$ws.Range("I102").Value = 17605
$ws.Range("K102").Value = 17605
$ws.Range("M102").Value = -15983
$ws.Range("H123").Value = 51302.273
$ws.Range("J123").Value = 51302.273
$ws.Range("L123").Value = 51302.273
$ws.Range("N123").Value = -56202.273
$ws.Range("H132").Value = 2092.6943
$ws.Range("I132").Value = 1457.2963
$ws.Range("K132").Value = 4371.8889
$ws.Range("M132").Value = -1841.8889

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2381.56
$ws.Range("I7").Value = 2416.6667
$ws.Range("K7").Value = 2416.6667
$ws.Range("M7").Value = -2304.6667
$ws.Range("H68").Value = 4589.1
$ws.Range("I68").Value = 3499
$ws.Range("J68").Value = 5056.2856
$ws.Range("K68").Value = 3499
$ws.Range("L68").Value = 5056.2856
$ws.Range("M68").Value = -2750
$ws.Range("N68").Value = -6554.2856
$ws.Range("H71").Value = 4589.1
$ws.Range("I71").Value = 3499
$ws.Range("J71").Value = 5056.2856
$ws.Range("K71").Value = 17495
$ws.Range("L71").Value = 25281.428
$ws.Range("M71").Value = -13751
$ws.Range("N71").Value = -32769.428
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H93").Value = 1272.6129
$ws.Range("I93").Value = 981
$ws.Range("K93").Value = 981
$ws.Range("M93").Value = 267
$ws.Range("H122").Value = 4114.7646
$ws.Range("I122").Value = 3381.6155
$ws.Range("J122").Value = 6497.5
$ws.Range("K122").Value = 10144.8465
$ws.Range("L122").Value = 19492.5
$ws.Range("M122").Value = -7694.8465
$ws.Range("N122").Value = -24392.5
$ws.Range("H126").Value = 2381.56
$ws.Range("I126").Value = 2416.6667
$ws.Range("K126").Value = 7250.000100000001
$ws.Range("M126").Value = -4780.000100000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H62").Value = 2987.5
$ws.Range("I62").Value = 2987.5
$ws.Range("K62").Value = 2987.5
$ws.Range("M62").Value = -2363.5
$ws.Range("H65").Value = 2987.5
$ws.Range("I65").Value = 2987.5
$ws.Range("K65").Value = 14937.5
$ws.Range("M65").Value = -11817.5
$ws.Range("H107").Value = 662.5294
$ws.Range("I107").Value = 466.76923
$ws.Range("K107").Value = 1400.30769
$ws.Range("M107").Value = 519.6923099999999
$ws.Range("H122").Value = 38154.227
$ws.Range("I122").Value = 46803.03
$ws.Range("J122").Value = 3559
$ws.Range("K122").Value = 140409.09
$ws.Range("L122").Value = 10677
$ws.Range("M122").Value = -137959.09
$ws.Range("N122").Value = -15577
$ws.Range("H126").Value = 230566.5
$ws.Range("I126").Value = 3025.3572
$ws.Range("K126").Value = 9076.071599999999
$ws.Range("M126").Value = -6606.071599999999
